# Updates cryptos list data (prices / volume changes), per commit
# "Updated cryptos list on Wed Mar  6 10:51:00 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.612.17"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "3.837.01"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "425.41"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.92"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "3.829.27"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.613"
$ws.Range("E8").Value = "  -5.02%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  -5.65%  "
$ws.Range("E11").Value = "  -8.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000370"
$ws.Range("E12").Value = "  -9.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.88"
$ws.Range("E13").Value = "  -4.84%  "
$ws.Range("D14").Value = "4.432.54"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.09"
$ws.Range("E15").Value = "  -5.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.63"
$ws.Range("E16").Value = "  +16.77%  "
$ws.Range("D17").Value = "3.841.28"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.59"
$ws.Range("E19").Value = "  -5.66%  "
$ws.Range("D20").Value = "66.867.47"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("E21").Value = "  -6.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "410.57"
$ws.Range("E22").Value = "  -8.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.47"
$ws.Range("E23").Value = "  -12.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.41"
$ws.Range("E24").Value = "  -4.93%  "
$ws.Range("E25").Value = "  -3.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.94"
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.66"
$ws.Range("E27").Value = "  +12.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.24"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("E29").Value = "  -6.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "688.94"
$ws.Range("E30").Value = "  +5.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.46"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.74"
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.20"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  -8.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.58"
$ws.Range("E36").Value = "  -8.16%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0802"
$ws.Range("E37").Value = "  +7.18%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.06"
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.14"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0457"
$ws.Range("E41").Value = "  -8.18%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  -8.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.94"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.50"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("E46").Value = "  -4.73%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.08"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.11"
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.33"
$ws.Range("E49").Value = "  -10.59%  "
$ws.Range("E50").Value = "  -4.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.54"
$ws.Range("E51").Value = "  -5.08%  "
